$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = "[-, 'MEC-1B-Comandos Eletricos', -, -]"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "['MEC-1B-Comandos Eletricos', -, -, 'MEC-2B-Elet. Dig. Bas.']"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "[-, -, 'MCT-3A-Robótica', -]"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "['MEC-1B-Comandos Eletricos', -, -, 'MEC-2B-Elet. Dig. Bas.']"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "[-, -, 'MCT-3A-Robótica', -]"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "[-, -, -, 'MEC-2B-Elet. Dig. Bas.']"
$ws.Range("F6").Value = "[-, -, 'MCT-3A-Robótica', -]"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "[-, -, -, 'MEC-2B-Elet. Dig. Bas.']"
$ws.Range("F7").Value = "[-, -, 'MCT-3A-Robótica', -]"

# Row 8
$ws.Range("F8").Value = "[-, -, 'MEC-1B-Comandos Eletricos', -]"

# Row 18
$ws.Range("B18").Value = "['ELM-2NA-Eletrônica Básica', 'ELM-2NA-Eletrônica Básica']"
$ws.Range("E18").Value = "[-, 'MEC-1NA-Comandos Eletricos', -, -]"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("B19").Value = "['ELM-2NA-Eletrônica Básica', 'ELM-2NA-Eletrônica Básica']"
$ws.Range("E19").Value = "[-, 'MEC-1NA-Comandos Eletricos', -, -]"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "[-, 'MEC-1NA-Comandos Eletricos', -, -]"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "[-, 'MEC-1NA-Comandos Eletricos', -, -]"
$ws.Range("F21").Value = "-"
